# feat: add 2022-Q4 data
#
# Starting layout:
#   Sheet1 = "总计"      (summary)
#   Sheet2 = "2022-Q1"   (fund holdings for 2022-Q1)
#
# Target layout:
#   Sheet1 = "总计"      (summary, gains a 2022-Q4 row above the 2022-Q1 row)
#   Sheet2 = "2022-Q4"   (new fund holdings table for 2022-Q4)
#   Sheet3 = "2022-Q1"   (original fund holdings table, unchanged, moved here)

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q1 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet right after the existing "2022-Q1" sheet.
#    It will inherit the NEXT sheetId (3) and become the new "2022-Q1" home,
#    keeping today's fund-holdings data exactly as it is.
# ---------------------------------------------------------------------------
$q1Copy = $wb.Worksheets.Add($null, $q1)

# Copy every used cell (values, then formatting) from the original sheet onto
# the freshly added one. (PasteSpecial(-4104) carries values/types;
# PasteSpecial(-4122) carries formats - doing both nets a full copy.)
$srcRange = $q1.Range("A1:H3")
$dstRange = $q1Copy.Range("A1:H3")
$srcRange.Copy()
$dstRange.PasteSpecial(-4104)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Re-purpose the original "2022-Q1" sheet (sheetId stays 2) as "2022-Q4"
#    and replace its contents with the new quarter's fund-holdings table.
#    (Rename the original sheet first so the "2022-Q1" name is free for the
#    freshly-added copy.)
# ---------------------------------------------------------------------------
$q4 = $q1
$q4.Name = "2022-Q4"
$q1Copy.Name = "2022-Q1"

$q4.Range("A1:H3").Clear()

$headerStyleSrc = $summary.Range("B1")

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$cols = @("B","C","D","E","F","G","H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
}
$q4.Range("B1:H1").Style = $headerStyleSrc.Style

$q4.Range("A2").Value = 0
$q4.Range("A2").Style = $summary.Range("A2").Style

$q4.Range("B2").Value = "001614"
$q4.Range("C2").Value = "东方区域发展混合"

$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.16"
$q4.Range("E2").Value = "93.97"
$q4.Range("F2").Value = "4.57"
$q4.Range("G2").Value = "0.0073"
$q4.Range("D2:G2").ClearFormats()

$q4.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: push the existing 2022-Q1 row down to
#    row 3, and add a new 2022-Q4 row in row 2.
# ---------------------------------------------------------------------------
$rowStyleSrc = $summary.Range("A2")

$summary.Range("A3").Value = 1
$summary.Range("A3").Style = $rowStyleSrc.Style
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

$wb.Worksheets.Item(1).Select()
